$wb = $excel.ActiveWorkbook

# Duplicate "Iteration 4" to create the new "Iteration 5" sheet, placed
# right before "Guidelines for Bug Metrics".
$src = $wb.Worksheets.Item("Iteration 4")
$guidelines = $wb.Worksheets.Item("Guidelines for Bug Metrics")
$src.Copy($guidelines)
$new = $wb.Worksheets.Item("Iteration 4 (2)")
$new.Name = "Iteration 5"

# Update the iteration title banner.
$new.Range("B2").Value = "Iteration 5 (15 December 2014 - 28 December 2014)"

# Replace the test case row with the new iteration's bug entry.
$new.Range("C8").Value = "State transition"
$new.Range("D8").Value = "editStateLecturer.jsp"
$new.Range("E8").Value = "Null pointer found when no scenario is activated"
$new.Range("F8").Value = "12/25/2014"
$new.Range("G8").Value = 5
$new.Range("I8").Value = "Unsolved"
$new.Range("J8").ClearContents()
$new.Range("K8").ClearContents()
$new.Range("L8").ClearContents()

$new.Rows(8).RowHeight = 30.75

# Move tab selection to the newly created iteration.
$src.Range("H8").Select()
$new.Activate()
$new.Range("E9").Select()
